$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price (D) column hold numeric-looking text (e.g. "174.80", "1.00").
# Force them to stay as Text before assignment, then clear the temporary
# number-format override so no extra style is left behind on the cell.
$priceCells = @("D2", "D3", "D4", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D48", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.899.43"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").Value = "3.319.58"
$ws.Range("E3").Value = "  -1.88%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("E5").Value = "  -1.74%  "

$ws.Range("D6").Value = "174.80"
$ws.Range("E6").Value = "  -6.99%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "0.582"
$ws.Range("E8").Value = "  -2.50%  "

$ws.Range("D9").Value = "3.312.25"
$ws.Range("E9").Value = "  -1.82%  "

$ws.Range("E10").Value = "  -4.98%  "

$ws.Range("D11").Value = "0.577"
$ws.Range("E11").Value = "  -2.42%  "

$ws.Range("D12").Value = "45.40"
$ws.Range("E12").Value = "  -4.95%  "

$ws.Range("D13").Value = "0.0000269"
$ws.Range("E13").Value = "  -3.19%  "

$ws.Range("D14").Value = "663.17"
$ws.Range("E14").Value = "  +3.57%  "

$ws.Range("D15").Value = "3.856.02"
$ws.Range("E15").Value = "  -1.83%  "

$ws.Range("E16").Value = "  -2.95%  "

$ws.Range("D17").Value = "68.004.32"
$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("E18").Value = "  -0.92%  "

$ws.Range("D19").Value = "3.318.54"
$ws.Range("E19").Value = "  -1.61%  "

$ws.Range("D20").Value = "17.48"
$ws.Range("E20").Value = "  -3.51%  "

$ws.Range("D21").Value = "10.87"
$ws.Range("E21").Value = "  -2.74%  "

$ws.Range("D22").Value = "0.888"
$ws.Range("E22").Value = "  -2.90%  "

$ws.Range("D23").Value = "17.11"
$ws.Range("E23").Value = "  -5.38%  "

$ws.Range("D24").Value = "5.35"
$ws.Range("E24").Value = "  +4.24%  "

$ws.Range("D25").Value = "97.44"
$ws.Range("E25").Value = "  -2.53%  "

$ws.Range("D26").Value = "3.84"
$ws.Range("E26").Value = "  -5.08%  "

$ws.Range("E27").Value = "  -6.41%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "33.87"
$ws.Range("E28").Value = "  +3.64%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "9.26"
$ws.Range("E29").Value = "  -5.85%  "

$ws.Range("D30").Value = "8.42"
$ws.Range("E30").Value = "  -3.84%  "

$ws.Range("D31").Value = "7.26"
$ws.Range("E31").Value = "  +3.87%  "

$ws.Range("D32").Value = "588.19"
$ws.Range("E32").Value = "  -4.27%  "

$ws.Range("D33").Value = "10.95"
$ws.Range("E33").Value = "  -1.80%  "

$ws.Range("E34").Value = "  -2.47%  "

$ws.Range("D35").Value = "3.732.18"
$ws.Range("E35").Value = "  -7.12%  "

$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "55.70"
$ws.Range("E37").Value = "  -0.92%  "

$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "3.34"
$ws.Range("E38").Value = "  -13.51%  "

$ws.Range("D39").Value = "0.132"
$ws.Range("E39").Value = "  -0.36%  "

$ws.Range("D40").Value = "32.57"
$ws.Range("E40").Value = "  -4.11%  "

$ws.Range("E41").Value = "  -7.02%  "

$ws.Range("D42").Value = "3.09"
$ws.Range("E42").Value = "  -5.22%  "

$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "0.333"
$ws.Range("E43").Value = "  -3.67%  "

$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0665"
$ws.Range("E44").Value = "  -6.49%  "

$ws.Range("E45").Value = "  -4.84%  "

$ws.Range("E46").Value = "  -4.28%  "

$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("D48").Value = "0.127"
$ws.Range("E48").Value = "  -2.55%  "

$ws.Range("E49").Value = "  +0.31%  "

$ws.Range("E50").Value = "  -2.28%  "

$ws.Range("D51").Value = "128.34"
$ws.Range("E51").Value = "  +0.04%  "

foreach ($addr in $priceCells) {
    $ws.Range($addr).ClearFormats()
}
